# Update the "Last Updated" timestamp on the Metadata sheet.
$wb = $excel.ActiveWorkbook
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("A2").Value = "05 Nov 2025, 02:08 PM"

# "Stock List" sheet: a new top entry (CAPTRU-RE1) was added, pushing every
# existing row down by one and dropping the final (77th) row that fell off
# the bottom of the list.
$ws = $wb.Worksheets.Item("Stock List")

# Insert a new blank row at row 2, shifting rows 2-76 down to 3-77.
$ws.Rows.Item(2).Insert()

# The old row 76 (TRAVELFOOD) is now duplicated at row 77; drop it so the
# sheet still ends at row 76.
$ws.Rows.Item(77).Delete()

# Populate the new row 2 with the new top entry's data.
$ws.Range("A2").Value = "📋"
$ws.Range("B2").Value = "CAPTRU-RE1"
$ws.Range("C2").Value = "CAPTRU-RE1"
$ws.Range("D2").Value = 5.67
$ws.Range("E2").Value = -11.9565
$ws.Range("F2").Value = "N/A"
$ws.Range("G2").Value = "N/A"
$ws.Range("H2").Value = 0
